# Scheduled-runner style refresh of market-price / leve-profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H..N)
# for the rows whose underlying market data changed since the last run.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3699.4375
$ws.Range("J17").Value = 3004.6428
$ws.Range("L17").Value = 9013.928400000001
$ws.Range("N17").Value = -9349.928400000001
$ws.Range("H33").Value = 188.1
$ws.Range("I33").Value = 188.1
$ws.Range("K33").Value = 188.1
$ws.Range("M33").Value = 40.90000000000001
$ws.Range("H88").Value = 2237.125
$ws.Range("I88").Value = 1499.5
$ws.Range("K88").Value = 1499.5
$ws.Range("M88").Value = -1093.5
$ws.Range("H91").Value = 2237.125
$ws.Range("I91").Value = 1499.5
$ws.Range("K91").Value = 1499.5
$ws.Range("M91").Value = -95.5
$ws.Range("H94").Value = 2497.6
$ws.Range("I94").Value = 2497.6
$ws.Range("K94").Value = 2497.6
$ws.Range("M94").Value = -2046.6
$ws.Range("H98").Value = 1371.1538
$ws.Range("I98").Value = 1160.238
$ws.Range("K98").Value = 1160.238
$ws.Range("M98").Value = 337.7619999999999
$ws.Range("H107").Value = 1145.7
$ws.Range("I107").Value = 744.625
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 744.625
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = 1175.375
$ws.Range("N107").Value = -6590
$ws.Range("H116").Value = 17785.572
$ws.Range("J116").Value = 4799.8
$ws.Range("L116").Value = 4799.8
$ws.Range("N116").Value = -11683.8
$ws.Range("H122").Value = 1371.1538
$ws.Range("I122").Value = 1160.238
$ws.Range("K122").Value = 3480.714
$ws.Range("M122").Value = -1030.714
$ws.Range("H129").Value = 1083.4546
$ws.Range("J129").Value = 1137.9487
$ws.Range("L129").Value = 3413.8461
$ws.Range("N129").Value = -13413.8461
$ws.Range("H132").Value = 997.24243
$ws.Range("I132").Value = 892.1923
$ws.Range("K132").Value = 2676.5769
$ws.Range("M132").Value = -146.5769
$ws.Range("H133").Value = 59480
$ws.Range("J133").Value = 59480
$ws.Range("L133").Value = 59480
$ws.Range("N133").Value = -69600
$ws.Range("H137").Value = 3054.818
$ws.Range("I137").Value = 2762.875
$ws.Range("K137").Value = 8288.625
$ws.Range("M137").Value = -5738.625
$ws.Range("H141").Value = 2804321
$ws.Range("I141").Value = 3501623.5
$ws.Range("K141").Value = 10504870.5
$ws.Range("M141").Value = -10499690.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 371769.6
$ws.Range("J2").Value = 1446.8334
$ws.Range("L2").Value = 1446.8334
$ws.Range("N2").Value = -1672.8334
$ws.Range("H32").Value = 4193.2456
$ws.Range("I32").Value = 3178.2593
$ws.Range("K32").Value = 3178.2593
$ws.Range("M32").Value = -2891.2593
$ws.Range("H74").Value = 1181.4333
$ws.Range("I74").Value = 474.86365
$ws.Range("J74").Value = 3124.5
$ws.Range("K74").Value = 474.86365
$ws.Range("L74").Value = 3124.5
$ws.Range("M74").Value = 399.13635
$ws.Range("N74").Value = -4872.5
$ws.Range("H77").Value = 1181.4333
$ws.Range("I77").Value = 474.86365
$ws.Range("J77").Value = 3124.5
$ws.Range("K77").Value = 2374.31825
$ws.Range("L77").Value = 15622.5
$ws.Range("M77").Value = 1993.68175
$ws.Range("N77").Value = -24358.5
$ws.Range("H109").Value = 58656.5
$ws.Range("J109").Value = 58656.5
$ws.Range("L109").Value = 58656.5
$ws.Range("N109").Value = -61430.5
$ws.Range("H110").Value = 296.27274
$ws.Range("I110").Value = 270.9
$ws.Range("J110").Value = 550
$ws.Range("K110").Value = 270.9
$ws.Range("L110").Value = 550
$ws.Range("M110").Value = 1774.1
$ws.Range("N110").Value = -4640
$ws.Range("H116").Value = 371769.6
$ws.Range("J116").Value = 1446.8334
$ws.Range("L116").Value = 1446.8334
$ws.Range("N116").Value = -6034.8334
$ws.Range("H122").Value = 785.5
$ws.Range("I122").Value = 599.48
$ws.Range("J122").Value = 1715.6
$ws.Range("K122").Value = 1798.44
$ws.Range("L122").Value = 5146.799999999999
$ws.Range("M122").Value = 651.5599999999999
$ws.Range("N122").Value = -10046.8
$ws.Range("H123").Value = 73500
$ws.Range("J123").Value = 73500
$ws.Range("L123").Value = 73500
$ws.Range("N123").Value = -83300
$ws.Range("H132").Value = 2160.762
$ws.Range("I132").Value = 1691
$ws.Range("J132").Value = 3335.1667
$ws.Range("K132").Value = 5073
$ws.Range("L132").Value = 10005.5001
$ws.Range("M132").Value = -2543
$ws.Range("N132").Value = -15065.5001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 371769.6
$ws.Range("J3").Value = 1446.8334
$ws.Range("L3").Value = 1446.8334
$ws.Range("N3").Value = -1674.8334
$ws.Range("H99").Value = 1435.3
$ws.Range("I99").Value = 1187.5
$ws.Range("J99").Value = 1600.5
$ws.Range("K99").Value = 1187.5
$ws.Range("L99").Value = 1600.5
$ws.Range("M99").Value = 310.5
$ws.Range("N99").Value = -4596.5
$ws.Range("H105").Value = 2224.9644
$ws.Range("I105").Value = 2039.0869
$ws.Range("K105").Value = 2039.0869
$ws.Range("M105").Value = -292.0869

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1133.7693
$ws.Range("I107").Value = 1005.7143
$ws.Range("K107").Value = 1005.7143
$ws.Range("M107").Value = 914.2857
$ws.Range("H132").Value = 2489.2917
$ws.Range("I132").Value = 1631.4706
$ws.Range("K132").Value = 4894.4118
$ws.Range("M132").Value = -2364.4118
$ws.Range("H141").Value = 63799
$ws.Range("J141").Value = 61748.75
$ws.Range("L141").Value = 61748.75
$ws.Range("N141").Value = -72108.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 946.25
$ws.Range("I122").Value = 486.5
$ws.Range("K122").Value = 4378.5
$ws.Range("M122").Value = -1928.5
$ws.Range("H131").Value = 13651.419
$ws.Range("I131").Value = 590
$ws.Range("J131").Value = 14552.207
$ws.Range("K131").Value = 1770
$ws.Range("L131").Value = 43656.621
$ws.Range("M131").Value = 3270
$ws.Range("N131").Value = -53736.621

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3500000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 4106352.8
$ws.Range("J7").Value = 1061599.6
$ws.Range("L7").Value = 1061599.6
$ws.Range("N7").Value = -1061823.6
$ws.Range("H8").Value = 4106352.8
$ws.Range("J8").Value = 1061599.6
$ws.Range("L8").Value = 1061599.6
$ws.Range("N8").Value = -1061877.6
$ws.Range("H97").Value = 579.03125
$ws.Range("I97").Value = 576.9666999999999
$ws.Range("J97").Value = 610
$ws.Range("K97").Value = 576.9666999999999
$ws.Range("L97").Value = 610
$ws.Range("M97").Value = -80.96669999999995
$ws.Range("N97").Value = -1602
$ws.Range("H113").Value = 1390.5555
$ws.Range("I113").Value = 1158
$ws.Range("J113").Value = 1506.8334
$ws.Range("K113").Value = 1158
$ws.Range("L113").Value = 1506.8334
$ws.Range("M113").Value = 1012
$ws.Range("N113").Value = -5846.8334
$ws.Range("H122").Value = 1802.8889
$ws.Range("I122").Value = 1775.7693
$ws.Range("J122").Value = 1873.4
$ws.Range("K122").Value = 5327.3079
$ws.Range("L122").Value = 5620.200000000001
$ws.Range("M122").Value = -2877.3079
$ws.Range("N122").Value = -10520.2

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10988.706
$ws.Range("I40").Value = 11267.833
$ws.Range("J40").Value = 10318.8
$ws.Range("K40").Value = 11267.833
$ws.Range("L40").Value = 10318.8
$ws.Range("M40").Value = -11131.833
$ws.Range("N40").Value = -10590.8
$ws.Range("H122").Value = 9560.789000000001
$ws.Range("I122").Value = 9040.9375
$ws.Range("K122").Value = 27122.8125
$ws.Range("M122").Value = -24672.8125
$ws.Range("H132").Value = 1575.5416
$ws.Range("I132").Value = 1244.9524
$ws.Range("K132").Value = 3734.857199999999
$ws.Range("M132").Value = -1204.857199999999
$ws.Range("H136").Value = 3242.1428
$ws.Range("I136").Value = 3358.2
$ws.Range("K136").Value = 10074.6
$ws.Range("M136").Value = -7524.599999999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 200062540
$ws.Range("J47").Value = 200062540
$ws.Range("L47").Value = 200062540
$ws.Range("N47").Value = -200063684
$ws.Range("H81").Value = 1149.3334
$ws.Range("I81").Value = 1499
$ws.Range("K81").Value = 2998
$ws.Range("M81").Value = -1937
$ws.Range("H84").Value = 1149.3334
$ws.Range("I84").Value = 1499
$ws.Range("K84").Value = 14990
$ws.Range("M84").Value = -9686
$ws.Range("H123").Value = 48084
$ws.Range("J123").Value = 48084
$ws.Range("L123").Value = 48084
$ws.Range("N123").Value = -57884
$ws.Range("H132").Value = 1288.3103
$ws.Range("I132").Value = 917.6512
$ws.Range("J132").Value = 2350.8667
$ws.Range("K132").Value = 2752.9536
$ws.Range("L132").Value = 7052.6001
$ws.Range("M132").Value = -222.9535999999998
$ws.Range("N132").Value = -12112.6001
